$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0596774193548387
$ws.Range("C2").Value = -0.151915455746367
$ws.Range("D2").Value = 0.204419889502762
$ws.Range("E2").Value = 0.125776397515528
$ws.Range("F2").Value = -0.121831252364737

$ws.Range("B3").Value = 0.246774193548387
$ws.Range("C3").Value = 0.264200792602378
$ws.Range("D3").Value = 0.558011049723757
$ws.Range("E3").Value = 0.301242236024845
$ws.Range("F3").Value = 0.0620506999621642

$ws.Range("B4").Value = 0.154838709677419
$ws.Range("C4").Value = 0.408190224570674
$ws.Range("D4").Value = 0.548802946593002
$ws.Range("E4").Value = 0.607142857142857
$ws.Range("F4").Value = 0.178584941354521

$ws.Range("B5").Value = 0.479032258064516
$ws.Range("C5").Value = 0.513870541611625
$ws.Range("D5").Value = 0.710865561694291
$ws.Range("E5").Value = 0.515527950310559
$ws.Range("F5").Value = 0.195611048051457

$ws.Range("B6").Value = 1.00806451612903
$ws.Range("C6").Value = 0.895640686922061
$ws.Range("D6").Value = 1.11786372007366
$ws.Range("E6").Value = 0.992236024844721
$ws.Range("F6").Value = 0.184638668180098

$ws.Range("B7").Value = 0.141935483870968
$ws.Range("C7").Value = 0.235138705416116
$ws.Range("D7").Value = 0.270718232044199
$ws.Range("E7").Value = 0.411490683229814
$ws.Range("F7").Value = -0.11312902005297
